$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45943
$ws.Range("B2").Value = 101.01
$ws.Range("C2").Value = 92.19
$ws.Range("D2").Value = 83.88
$ws.Range("E2").Value = 92.17
$ws.Range("F2").Value = 95.06999999999999
$ws.Range("G2").Value = 101.93
$ws.Range("H2").Value = 105.92
$ws.Range("I2").Value = 115.2
$ws.Range("J2").Value = 127.41
$ws.Range("K2").Value = 120.55
$ws.Range("L2").Value = 103.94
$ws.Range("M2").Value = 88.36
$ws.Range("N2").Value = 77.29000000000001
$ws.Range("O2").Value = 70.42
$ws.Range("P2").Value = 68.7
$ws.Range("Q2").Value = 65.47
$ws.Range("R2").Value = 71.14
$ws.Range("S2").Value = 84.3
$ws.Range("T2").Value = 99.91
$ws.Range("U2").Value = 132.08
$ws.Range("V2").Value = 162.09
$ws.Range("W2").Value = 160.49
$ws.Range("X2").Value = 129.25
$ws.Range("Y2").Value = 115.02
$ws.Range("Z2").Value = 102.66
$ws.Range("AB2").Value = 141.71
$ws.Range("AD2").Value = 161.29
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 123.98
$ws.Range("AG2").Value = "0h-18h"

$wb.Save()
